# Update cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values like
# "1.00" or "10.60" are not coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$rows = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="61.935.53"; E="  -8.75%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="3.173.83"; E="  -10.28%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.00"; E="  -0.28%  "},
    @{Row=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="505.81"; E="  -8.99%  "},
    @{Row=6; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="169.06"; E="  -14.28%  "},
    @{Row=7; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.582"; E="  -11.55%  "},
    @{Row=8; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.00"; E="  -0.08%  "},
    @{Row=9; B="LidoStakedEther"; C="https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D="3.175.21"; E="  -10.10%  "},
    @{Row=10; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.585"; E="  -11.79%  "},
    @{Row=11; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="53.69"; E="  -11.67%  "},
    @{Row=12; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.128"; E="  -11.35%  "},
    @{Row=13; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.0000248"; E="  -8.21%  "},
    @{Row=14; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="8.74"; E="  -12.12%  "},
    @{Row=15; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="3.669.37"; E="  -10.67%  "},
    @{Row=16; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="3.159.59"; E="  -10.85%  "},
    @{Row=17; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.112"; E="  -10.18%  "},
    @{Row=18; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="61.810.00"; E="  -8.75%  "},
    @{Row=19; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="16.74"; E="  -9.10%  "},
    @{Row=20; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="10.62"; E="  -10.91%  "},
    @{Row=21; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.925"; E="  -10.54%  "},
    @{Row=22; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="355.69"; E="  -11.07%  "},
    @{Row=23; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="3.61"; E="  -9.99%  "},
    @{Row=24; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="78.04"; E="  -10.00%  "},
    @{Row=25; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="6.07"; E="  -1.67%  "},
    @{Row=26; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="10.60"; E="  -10.00%  "},
    @{Row=27; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="3.78"; E="  -2.31%  "},
    @{Row=28; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="2.56"; E="  -9.88%  "},
    @{Row=29; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="10.89"; E="  -12.48%  "},
    @{Row=30; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="7.96"; E="  -10.90%  "},
    @{Row=31; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="27.55"; E="  -12.16%  "},
    @{Row=32; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="606.69"; E="  -16.06%  "},
    @{Row=33; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="6.30"; E="  -10.89%  "},
    @{Row=34; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="10.83"; E="  -8.09%  "},
    @{Row=35; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.00"; E="  +0.06%  "},
    @{Row=36; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="56.34"; E="  -12.49%  "},
    @{Row=37; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.101"; E="  -10.19%  "},
    @{Row=38; B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="35.66"; E="  -7.86%  "},
    @{Row=39; B="TheGraph"; C="https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"; D="0.369"; E="  -6.14%  "},
    @{Row=40; B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="0.996"; E="  -0.24%  "},
    @{Row=41; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="0.0₃0655"; E="  -4.59%  "},
    @{Row=42; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.118"; E="  -10.85%  "},
    @{Row=43; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="2.777.81"; E="  -9.84%  "},
    @{Row=44; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="2.36"; E="  -5.92%  "},
    @{Row=45; B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="2.56"; E="  -7.41%  "},
    @{Row=46; B="ThetaToken"; C="https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"; D="2.55"; E="  -15.68%  "},
    @{Row=47; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0377"; E="  -8.20%  "},
    @{Row=48; B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="2.68"; E="  +0.55%  "},
    @{Row=49; B="ApeXProtocol"; C="https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; D="2.87"; E="  -5.52%  "},
    @{Row=50; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="131.77"; E="  -5.57%  "},
    @{Row=51; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.120"; E="  -11.83%  "}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
